# "fix button col width": the "Copy" (row 19) and "Delete" (row 20) button
# rows on the "Scouting Admin" sheet were missing their "Works" date in
# column B, unlike every other feature row. Backfill them with the same
# date (1/26/2024 -> serial 45317) used by the surrounding rows, copying
# the existing date formatting from B3 so the new cells reuse the same
# style record instead of creating a new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scouting Admin")

# Copy the number-formatted style already used by the other date cells
# (e.g. B3) onto B19:B20 so they match the rest of the column.
$ws.Range("B3").Copy()
$ws.Range("B19:B20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the missing "Works" dates.
$ws.Range("B19").Value = 45317
$ws.Range("B20").Value = 45317
